$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new data row (row 86) with titration results for blue tank CRM opened 20220427 dmbp
$ws.Range("A86").Value2 = 20220608
$ws.Range("B86").Value2 = 2228.8510392501998
$ws.Range("C86").Value2 = 2224.4699999999998
$ws.Range("D86").Formula = "=100*(B86-C86)/C86"
$ws.Range("E86").Value2 = 180
$ws.Range("F86").Value2 = "CRM OPENED 20220427 dmbp"

# Update column A width (autofit-style widening seen in the saved workbook)
$ws.Columns.Item(1).ColumnWidth = 9.14

# Update the window selection / view state to match the post-edit session
$ws.Range("C89").Select()
